$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of issue data appended below the existing table (row 14).
# Column A ("18") looks numeric, so prefix it to force text entry,
# then clear the resulting "quote prefix" style so no new cell style
# is introduced (matches the rest of the sheet, which has no explicit
# per-cell styling).
$ws.Range("A14").Value = "'18"
$ws.Range("A14").Style = "Normal"

$ws.Range("B14").Value = "hjasgdjahgsdhjags"
$ws.Range("C14").Value = "open"
$ws.Range("D14").Value = "2025-03-25T11:48:19Z"
$ws.Range("E14").Value = "bug"
